$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("exclusion_rules")

# New column E: "Ignore Store Policy" flag.
$ws.Range("E1").Value = "Ignore Store Policy"

# Rows 2-12 and 17-21 are "Include" rules -> Ignore Store Policy = 1
# Rows 13-16 are "Exclude" brand_name rules -> Ignore Store Policy = 0
$includeRows = @(2,3,4,5,6,7,8,9,10,11,12,17,18,19,20,21)
foreach ($r in $includeRows) {
    $ws.Cells.Item($r, 5).Value = 1
}

$excludeRows = @(13,14,15,16)
foreach ($r in $excludeRows) {
    $ws.Cells.Item($r, 5).Value = 0
}

# Update the active selection to reflect the last edited cell.
$ws.Range("E21").Select()

# store_policy sheet selection moves too.
$ws2 = $wb.Worksheets.Item("store_policy")
$ws2.Range("A3").Select()
$ws.Activate()
